$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Cells.Item(2, 6).Value = 162
    $ws.Cells.Item(3, 6).Value = 7340
    $ws.Cells.Item(4, 6).Value = 5696
    $ws.Cells.Item(5, 6).Value = 84
    $ws.Cells.Item(9, 6).Value = 113
    $ws.Cells.Item(11, 6).Value = 116
    $ws.Cells.Item(12, 6).Value = 207
    $ws.Cells.Item(13, 6).Value = 66
    $ws.Cells.Item(15, 6).Value = 390
    $ws.Cells.Item(20, 6).Value = 56
}
